# Applies the "Updated cryptos list" refresh: per-row Price (D) and
# Volume(1h) (E) updates, plus the Coin/Link (B/C) swap for rows 14-15
# (WrappedliquidstakedEther2.0 <-> Polkadot) as captured in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.125.58"
$ws.Range("E2").Value = "  -0.77%  "

# Row 3
$ws.Range("D3").Value = "3.767.50"
$ws.Range("E3").Value = "  +2.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.55"
$ws.Range("E5").Value = "  -2.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.06"
$ws.Range("E6").Value = "  +2.35%  "

# Row 7
$ws.Range("D7").Value = "3.758.57"
$ws.Range("E7").Value = "  +1.89%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  -4.97%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.730"
$ws.Range("E10").Value = "  -3.84%  "

# Row 11
$ws.Range("E11").Value = "  -8.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000359"
$ws.Range("E12").Value = "  -8.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.07"
$ws.Range("E13").Value = "  -3.70%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.02"
$ws.Range("E14").Value = "  -4.41%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.347.28"
$ws.Range("E15").Value = "  +1.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.75"
$ws.Range("E16").Value = "  +12.03%  "

# Row 17
$ws.Range("E17").Value = "  -1.30%  "

# Row 18
$ws.Range("D18").Value = "3.751.32"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.45"
$ws.Range("E19").Value = "  -4.95%  "

# Row 20
$ws.Range("D20").Value = "66.299.96"
$ws.Range("E20").Value = "  -0.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("E21").Value = "  -4.88%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "412.51"
$ws.Range("E22").Value = "  -7.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.42"
$ws.Range("E23").Value = "  -9.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.35"
$ws.Range("E24").Value = "  -4.82%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.07"
$ws.Range("E25").Value = "  -1.70%  "

# Row 26
$ws.Range("E26").Value = "  +14.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "35.88"
$ws.Range("E27").Value = "  -4.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.15"
$ws.Range("E28").Value = "  -4.91%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("E29").Value = "  -8.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "735.95"
$ws.Range("E30").Value = "  +12.68%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.36"
$ws.Range("E32").Value = "  -1.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.41"
$ws.Range("E34").Value = "  +1.87%  "

# Row 35
$ws.Range("E35").Value = "  -6.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.20"
$ws.Range("E36").Value = "  -3.93%  "

# Row 37
$ws.Range("E37").Value = "  +0.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.14"
$ws.Range("E38").Value = "  -3.91%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0735"
$ws.Range("E39").Value = "  -0.56%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0460"
$ws.Range("E40").Value = "  -6.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("E41").Value = "  -12.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -8.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.14"
$ws.Range("E44").Value = "  -7.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.08"
$ws.Range("E45").Value = "  -1.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  +19.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  -3.60%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("E48").Value = "  -1.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").Value = "  +1.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.26"
$ws.Range("E50").Value = "  -2.12%  "

# Row 51
$ws.Range("E51").Value = "  -3.33%  "
